$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Nacional" row (row 2) figures.
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 1.19
$ws.Range("E2").Value = 4435
$ws.Range("F2").Value = 35478
$ws.Range("G2").Value = 35478
$ws.Range("H2").Value = 44358.589647411864

# 2. Simplify the J3 total formula (same result, different formula text).
$ws.Range("J3").Formula = "=F3+F4"

# 3. Clear the stray empty formatted cells in column I (I2:I10).
$ws.Range("I2:I10").ClearFormats()

# 4. Highlight column F the same way column G is already highlighted,
#    row by row, so each row's F:G pair shares one fill color.
$ws.Range("F2").Interior.Color = $ws.Range("G2").Interior.Color
$ws.Range("F3").Interior.Color = $ws.Range("G3").Interior.Color
$ws.Range("F4").Interior.Color = $ws.Range("G4").Interior.Color
$ws.Range("F5:F10").Interior.Color = $ws.Range("G5").Interior.Color

# 5. J5's total now matches the centered look of the rest of the table.
$ws.Range("J5").HorizontalAlignment = -4108

# 6. Leave the selection on F2:G2, matching where the edit was made.
$ws.Range("F2:G2").Select()
